$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row: CornerCut (introduces "CornerCut" then "inch" shared strings) ---
$ws.Range("A15").Value = "CornerCut"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "inch"

# ScintillatorWidth / ScintillatorHeight: 145 mm -> 5.75 inch
$ws.Range("B5").Value = 5.75
$ws.Range("C5").Value = "inch"
$ws.Range("B6").Value = 5.75
$ws.Range("C6").Value = "inch"

# --- Fix label typos (shared-string text changes) ---
$ws.Range("A3").Value = "OpticalFiberClearance"
$ws.Range("A4").Value = "ScintillatorEdgeClearance"
$ws.Range("A12").Value = "FiberSensorClerance"

# --- Update values / units ---
# ScintillatorThickness: 7 -> 12.7 mm, with a new note in column D
$ws.Range("B2").Value = 12.7
$ws.Range("D2").Value = "(.5 in)"

# LoopLargeTrackRatio: 1.1 -> 1.25
$ws.Range("B10").Value = 1.25

# MountingScrewOffsett: 12.5 -> 9
$ws.Range("B13").Value = 9

# --- Selection moves to B4 ---
$ws.Range("B4").Select() | Out-Null

$wb.Save() | Out-Null
